$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.295.67"
$ws.Range("E2").Value = "  +0.44%  "

$ws.Range("D3").Value = "1.879.21"
$ws.Range("E3").Value = "  -1.20%  "

$ws.Range("E4").Value = "  -0.56%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.01"
$ws.Range("E5").Value = "  -3.05%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.679"
$ws.Range("E6").Value = "  -3.72%  "

$ws.Range("E7").Value = "  -0.56%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "43.79"
$ws.Range("E8").Value = "  +4.70%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.359"
$ws.Range("E9").Value = "  +1.01%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "53.52"
$ws.Range("E10").Value = "  +2.12%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0740"
$ws.Range("E11").Value = "  -2.81%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0976"
$ws.Range("E12").Value = "  -0.04%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "13.58"
$ws.Range("E13").Value = "  +3.03%  "

$ws.Range("D14").Value = "2.152.29"
$ws.Range("E14").Value = "  -1.26%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.766"
$ws.Range("E15").Value = "  +4.58%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.93"
$ws.Range("E16").Value = "  -1.10%  "

$ws.Range("D17").Value = "1.873.76"
$ws.Range("E17").Value = "  -1.54%  "

$ws.Range("D18").Value = "35.309.14"
$ws.Range("E18").Value = "  +0.42%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "73.35"
$ws.Range("E19").Value = "  -0.81%  "

$ws.Range("D20").Value = "0.0₃0822"
$ws.Range("E20").Value = "  -2.40%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "244.10"
$ws.Range("E21").Value = "  -0.52%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.84"
$ws.Range("E22").Value = "  -1.79%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.04"
$ws.Range("E23").Value = "  -0.03%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.68"
$ws.Range("E24").Value = "  +10.51%  "

$ws.Range("E25").Value = "  -0.49%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.16"
$ws.Range("E26").Value = "  -5.66%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "165.04"
$ws.Range("E27").Value = "  -2.28%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.59"
$ws.Range("E28").Value = "  +0.06%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.26"
$ws.Range("E29").Value = "  -1.41%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.127"
$ws.Range("E30").Value = "  -2.43%  "

$ws.Range("E31").Value = "  -1.38%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0591"
$ws.Range("E32").Value = "  -0.96%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.17"
$ws.Range("E33").Value = "  -2.62%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("B34").Value = "BinanceUSD"
$ws.Range("C34").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  -0.61%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("B35").Value = "WEMIXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").Value = "1.82"
$ws.Range("E35").Value = "  -13.72%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.42"
$ws.Range("E36").Value = "  -12.98%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.850"
$ws.Range("E37").Value = "  +0.35%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").Value = "1.94"
$ws.Range("E38").Value = "  -4.11%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").Value = "0.0729"
$ws.Range("E39").Value = "  +8.75%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.39"
$ws.Range("E40").Value = "  -0.61%  "

$ws.Range("E41").Value = "  +0.89%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "96.38"
$ws.Range("E42").Value = "  -2.35%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.08"
$ws.Range("E43").Value = "  -2.66%  "

$ws.Range("D44").Value = "1.306.07"
$ws.Range("E44").Value = "  -0.18%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.38"
$ws.Range("E45").Value = "  -1.15%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0797"
$ws.Range("E46").Value = "  +5.03%  "

$ws.Range("E47").Value = "  -1.00%  "

$ws.Range("E48").Value = "  -0.77%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "11.70"
$ws.Range("E49").Value = "  -3.73%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.27"
$ws.Range("E50").Value = "  -5.04%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "42.04"
$ws.Range("E51").Value = "  -1.98%  "
